# Update Saldo.xlsx "Export" sheet:
#  - Row 2 (LEILA):   Saldo 80039.86 -> 63000
#  - Reorder rows 4-6 from (ANA, ANDRE, GUSTAVO) to (ANDRE, GUSTAVO, ANA)
#    and update ANDRE's Saldo from 7726.01 -> 5026.01
#  - ANA's Saldo becomes 20000 -> 3000 (she now lands on what was row 6)
#  - GUSTAVO keeps his Saldo of 4292

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LEILA's balance changes.
$ws.Range("C2").Value = 63000

# Rows 4-6 get rewritten in their new order: ANDRE, GUSTAVO, ANA.
# Force text formatting on column A so the zero-padded account numbers
# don't get auto-coerced into plain numbers (losing their leading zeros).
$ws.Range("A4:A6").NumberFormat = "@"

$ws.Range("A4").Value = "005040864"
$ws.Range("B4").Value = "ANDRE"
$ws.Range("C4").Value = 5026.01

$ws.Range("A5").Value = "004313254"
$ws.Range("B5").Value = "GUSTAVO"
$ws.Range("C5").Value = 4292

$ws.Range("A6").Value = "005009922"
$ws.Range("B6").Value = "ANA"
$ws.Range("C6").Value = 3000
